$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 13, shifting existing rows 13-25 down to 15-27.
$ws.Rows("13:14").Insert()

# The insert operation leaves stray empty/bold-styled cells in column A for the two
# new rows (inherited from the row above); clear them since these rows have no label.
$ws.Range("A13:A14").Clear()

# The insert also mis-styles the new B/C cells (bold, no red/wrap) -- copy the correct
# number/wrap formats from row 10 (B = plain wrap, C = red wrap) down into rows 13-14.
$ws.Range("B10:C10").Copy()
$ws.Range("B13:C14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 10 (Objetivos:): fix B/C content with the full objectives paragraph.
$ws.Range("B10:C10").Value = 'Apresentar a análise química como ferramenta para o estudo da composição e das propriedades de materiais.Desenvolver a competência para formular e compreender problemas relacionados à análise química e buscar de forma autônoma procedimentos adequados para a sua solução. Desenvolver nos alunos a competência técnica para propor experimentos, obter e interpretar resultados analíticos. Incentivar trabalhos em grupo para a solução de problemas, com apresentação de resultados de forma oral e escrita.'

# New row 13: professor responsible for the course (Angelo).
$ws.Range("B13:C13").Value = '5840712 - Ângelo Capri Neto'

# New row 14: professor responsible for the course (Rosa).
$ws.Range("B14:C14").Value = '5840521 - Rosa Ana Conte'

# Row 15 (Programa resumido:): fix B/C content with the short-syllabus text.
$ws.Range("B15:C15").Value = '1. Introdução à Química Analítica;2. Preparação de amostras sólidas e líquidas;3. Métodos de análises qualitativas e quantitativas por via úmida;4. Métodos espectroscópicos de análise;5. Análise de gases em metais;'

# Row 17 (Programa:): fix B/C content with the full syllabus text.
$ws.Range("B17:C17").Value = 'Introdução à química analítica; Preparação de amostras sólidas e líquidas; Química analítica quantitativa por via úmida: Gravimetria e Volumetria; Métodos espectroscópicos de análise: interação radiação/matéria, absorção atômica e molecular. Espectroscopia UV/Visível: lei de Beer; instrumentação, calibração do equipamento, aplicações e interpretação dos resultados analíticos. Absorção Atômica: instrumentação, calibração do equipamento, identificação e controle de interferências; aplicações e interpretação de resultados analíticos. Emissão Atômica: instrumentação, calibração do equipamento e controle de interferências; aplicações e interpretação de resultados analíticos. Análise de gases em metais: instrumentação e calibração do equipamento; aplicações e interpretação de resultados analíticos.'

# Row 20 (Método:): fix B/C content with the evaluation method text.
$ws.Range("B20:C20").Value = 'A avaliação será feita por meio de duas provas (P1 e P2). A critério do professor, a avaliação poderá ser complementada por meio de trabalhos e/ou relatórios, valendo até 30% da nota das provas.'

# Row 21 (Critério:): fix B/C content with the final-grade criteria text.
$ws.Range("B21:C21").Value = 'A nota final (NF) será calculada pela média aritmética das provas. NF=(P1 +P2)/2.'

# Row 22 (Norma de recuperação:): add B/C content with the recovery-norm text.
$ws.Range("B22:C22").Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.'

# Row 23 (Bibliografia:): add B/C content with the bibliography text.
$ws.Range("B23:C23").Value = '1. VOGEL, A. L., et al. Análise Química Quantitativa, 6ª Ed., Rio de Janeiro, Livros Técnicos e Científicos Editora S.A., 2003.2. SKOOG, D.A. & Jeary, J.J. Principles of Instrumental Analysis, 6th Ed, Saunders College Publishing, 2007.3. MITRA, S. Sample Preparation Techniques in Analytical  Chemistry, Wiley & Sons, Hoboken, New Jersey, 2003.4. ANDERSON, R. Sample Pretreatment and  separation, Wiley & Sons, New York, 1997'

# Rows 25-27 (Requisitos list: LOB1012, LOQ4095, LOQ4098) already carry the correct
# text after the row insert shifted them down from rows 23-25, so no further edits
# are needed there.
